$wb = $excel.ActiveWorkbook

# Insert a new (blank) column before column N (14th) on the "Repayment schedule"
# sheet, shifting the old N/O/P ("Late"/"heading"/"Outstanding") columns right
# to O/P/Q.
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns.Item(14).Insert() | Out-Null
$wsRepay.Columns.Item(14).ColumnWidth = 9.83

# Re-point the active sheet from "Transactions" back to "Repayment schedule",
# keeping each sheet's own remembered selection (Transactions stays on C10,
# Repayment schedule ends up on S7).
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Select() | Out-Null
$wsTrans.Range("C10").Select() | Out-Null

$wsRepay.Select() | Out-Null
$wsRepay.Range("S7").Select() | Out-Null
